$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data (columns B through AC) between row 11 and row 12,
# leaving column A (the sequential index) untouched.
$r1 = $ws.Range("B11:AC11")
$r2 = $ws.Range("B12:AC12")
$vals1 = $r1.Value2
$vals2 = $r2.Value2
$r1.Value2 = $vals2
$r2.Value2 = $vals1

# Swap the data (columns B through AC) between row 83 and row 84,
# leaving column A (the sequential index) untouched.
$r3 = $ws.Range("B83:AC83")
$r4 = $ws.Range("B84:AC84")
$vals3 = $r3.Value2
$vals4 = $r4.Value2
$r3.Value2 = $vals4
$r4.Value2 = $vals3
